# Insert a new weekly price record for "Feria Lagunitas de Puerto Montt - Choclo"
# as row 138 (pushing the former rows 138..161 down to 139..162).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data down by inserting a blank row at position 138.
$ws.Rows.Item(138).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Cells.Item(138, 1).Value  = 4
$ws.Cells.Item(138, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(138, 3).Value  = "Los Lagos"
$ws.Cells.Item(138, 4).Value  = 44551
$ws.Cells.Item(138, 5).Value  = 10
$ws.Cells.Item(138, 6).Value  = 100112024
$ws.Cells.Item(138, 7).Value  = "Choclo"
$ws.Cells.Item(138, 8).Value  = "Dulce o Americano"
$ws.Cells.Item(138, 9).Value  = "Primera"
$ws.Cells.Item(138, 10).Value = 300
$ws.Cells.Item(138, 11).Value = 20000
$ws.Cells.Item(138, 12).Value = 20000
$ws.Cells.Item(138, 13).Value = 20000
$ws.Cells.Item(138, 14).Value = "$/malla 70 unidades"
$ws.Cells.Item(138, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(138, 16).Value = 286
$ws.Cells.Item(138, 17).Value = 70
$ws.Cells.Item(138, 18).Value = "Hortaliza"
